$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 and J1, matching the formatting of the existing header cells
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data values for columns I (I0) and J (IF), rows 2-41
$I0 = @(8,8,9,9,9,10,9,8,9,8,9,9,9,8,9,9,9,9,9,8,8,9,9,8,8,9,9,9,9,9,8,9,6,6,7,6,9,7,7,6)
$IF = @(9,8,10,9,9,11,9,8,9,9,9,9,9,9,9,9,9,9,9,9,8,9,9,8,8,9,9,9,9,9,8,9,7,6,8,6,9,7,7,7)

for ($i = 0; $i -lt 40; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $I0[$i]
    $ws.Cells.Item($row, 10).Value = $IF[$i]
}
